$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  "Jrue Holiday",       "PG,SG", "Boston Celtics"),
    @(3,  "Tyrese Haliburton",  "PG,SG", "Indiana Pacers"),
    @(4,  "Stephen Curry",      "PG,SG", "Golden State Warriors"),
    @(5,  "Darius Garland",     "PG",    "Cleveland Cavaliers"),
    @(6,  "OG Anunoby",         "SF,PF", "New York Knicks"),
    @(7,  "Daniel Gafford",     "PF,C",  "Dallas Mavericks"),
    @(8,  "Kevin Durant",       "SF,PF", "Phoenix Suns"),
    @(9,  "Jalen Johnson",      "SF,PF", "Atlanta Hawks"),
    @(10, "Jarrett Allen",      "C",     "Cleveland Cavaliers"),
    @(11, "Jalen Duren",        "C",     "Detroit Pistons"),
    @(12, "Trey Murphy III",    "SF,PF", "New Orleans Pelicans"),
    @(13, "Tyrese Maxey",       "PG,SG", "Philadelphia 76ers"),
    @(14, "Austin Reaves",      "PG,SG", "Los Angeles Lakers"),
    @(15, "Mark Williams",      "C",     "Charlotte Hornets"),
    @(16, "Keegan Murray",      "SF,PF", "Sacramento Kings"),
    @(17, "Franz Wagner",       "SF,PF", "Orlando Magic"),
    @(18, "Dereck Lively II",   "C",     "Dallas Mavericks"),
    @(19, "Karl-Anthony Towns", "PF,C",  "New York Knicks")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
